$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "17×13=",
    "61×79=",
    "49×77=",
    "40×45=",
    "13×30=",
    "51×65=",
    "38×75=",
    "97×41=",
    "71×80=",
    "91×79=",
    "62×77=",
    "50×70=",
    "31×54=",
    "84×97=",
    "28×24=",
    "20×75=",
    "94×55=",
    "55×53=",
    "44×66=",
    "49×53=",
    "11×49=",
    "66×38=",
    "32×52=",
    "28×77=",
    "37×81=",
    "88×28=",
    "35×95=",
    "41×17=",
    "59×26=",
    "69×62=",
    "23×87=",
    "56×36=",
    "19×25=",
    "85×42=",
    "19×31=",
    "86×71=",
    "87×32=",
    "79×35=",
    "43×57=",
    "80×38=",
    "66×38=",
    "35×27=",
    "73×28=",
    "40×45=",
    "31×67=",
    "92×93=",
    "30×19=",
    "58×90=",
    "29×38=",
    "62×94=",
    "11×47=",
    "21×15=",
    "42×97=",
    "91×19=",
    "50×100=",
    "78×72=",
    "38×64=",
    "93×73=",
    "56×41=",
    "76×79=",
    "45×53=",
    "12×22=",
    "28×79=",
    "73×80=",
    "96×16=",
    "61×11=",
    "18×12=",
    "35×14=",
    "63×13=",
    "92×90=",
    "22×48=",
    "51×95=",
    "99×31=",
    "48×74=",
    "31×76=",
    "61×94=",
    "92×76=",
    "35×92=",
    "98×67=",
    "34×76=",
    "27×54=",
    "78×78=",
    "90×20=",
    "61×96=",
    "41×56=",
    "75×95=",
    "82×78=",
    "11×56=",
    "45×83=",
    "59×62=",
    "64×88=",
    "52×28=",
    "63×59=",
    "50×86=",
    "52×90=",
    "19×47=",
    "95×11=",
    "25×68=",
    "40×95=",
    "75×46="
)

$rows = 20
$cols = 5
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $idx = ($r - 1) * $cols + ($c - 1)
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
    }
}

Write-Host "Done updating cells"
